# Apply crypto price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Preserve the existing "text" storage of the Price column (D) values:
# Excel auto-converts numeric-looking strings to real numbers when set via
# .Value, so force text format on the Price column first, then restore the
# original (default) style afterwards so no stray formatting is introduced.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.490.53"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "1.840.91"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "260.56"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").Value = "0.5253"
$ws.Range("E7").Value = "  +1.04%  "

$ws.Range("D8").Value = "0.3206"
$ws.Range("E8").Value = "  -0.42%  "

$ws.Range("D9").Value = "0.06791"
$ws.Range("E9").Value = "  +0.58%  "

$ws.Range("D10").Value = "18.80"
$ws.Range("E10").Value = "  +1.40%  "

$ws.Range("D11").Value = "0.7860"
$ws.Range("E11").Value = "  +3.29%  "

$ws.Range("D12").Value = "0.07749"
$ws.Range("E12").Value = "  +1.09%  "

$ws.Range("D13").Value = "1.844.44"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("D14").Value = "87.72"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").Value = "5.015"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").Value = "13.86"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("D19").Value = "0.000007944"
$ws.Range("E19").Value = "  +0.57%  "

$ws.Range("D20").Value = "26.512.52"
$ws.Range("E20").Value = "  +0.42%  "

$ws.Range("D21").Value = "2.072.64"
$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("D22").Value = "4.632"
$ws.Range("E22").Value = "  +1.96%  "

$ws.Range("D23").Value = "5.985"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("D24").Value = "9.384"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("D25").Value = "141.27"
$ws.Range("E25").Value = "  -2.18%  "

$ws.Range("D26").Value = "2.170"
$ws.Range("E26").Value = "  -2.89%  "

$ws.Range("D27").Value = "1.679"
$ws.Range("E27").Value = "  +1.70%  "

$ws.Range("D28").Value = "16.93"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "111.71"
$ws.Range("E29").Value = "  +0.50%  "

$ws.Range("D30").Value = "4.153"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").Value = "0.08686"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("D32").Value = "4.075"
$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").Value = "0.04868"
$ws.Range("E33").Value = "  +1.38%  "

$ws.Range("D34").Value = "0.7292"
$ws.Range("E34").Value = "  +4.33%  "

$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  +1.71%  "

$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("D37").Value = "3.090"
$ws.Range("E37").Value = "  +1.09%  "

$ws.Range("D38").Value = "2.249"
$ws.Range("E38").Value = "  +2.80%  "

$ws.Range("D39").Value = "0.01755"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").Value = "0.4774"
$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("E41").Value = "  +0.93%  "

$ws.Range("D42").Value = "109.55"
$ws.Range("E42").Value = "  -1.44%  "

$ws.Range("D43").Value = "5.936"
$ws.Range("E43").Value = "  -2.49%  "

$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").Value = "7.694"
$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("D46").Value = "0.4168"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "8.961"
$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05849"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.1232"
$ws.Range("E49").Value = "  +1.49%  "

$ws.Range("D50").Value = "34.87"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("D51").Value = "0.8936"
$ws.Range("E51").Value = "  +1.49%  "

# Restore default styling on the Price column now that the text values are set.
$priceRange.Style = "Normal"

